$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the east/north coordinates in row 5 to whole numbers
$ws.Range("Q5").Value = 528534
$ws.Range("R5").Value = 6936241

# Clear the start time (Z5) and end time (AB5) cells for row 5
$ws.Range("Z5").ClearContents()
$ws.Range("AB5").ClearContents()
